# Updated the page names for Subpopulation and LOT pages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LOT section (row 2, 11, 19, 26 - F/G/H columns) -> pop_filter2_section*
$ws.Range("F2").Value = "pop_filter2_section1"
$ws.Range("G2").Value = "pop_filter2_section1_checkbox"
$ws.Range("H2").Value = "pop_filter2_section"

$ws.Range("F11").Value = "pop_filter2_section1"
$ws.Range("G11").Value = "pop_filter2_section1_checkbox"
$ws.Range("H11").Value = "pop_filter2_section"

$ws.Range("F19").Value = "pop_filter2_section1"
$ws.Range("G19").Value = "pop_filter2_section1_checkbox"
$ws.Range("H19").Value = "pop_filter2_section"

$ws.Range("F26").Value = "pop_filter2_section1"
$ws.Range("G26").Value = "pop_filter2_section1_checkbox"
$ws.Range("H26").Value = "pop_filter2_section"

# Subpopulation section (row 10, 18 - F/G/H columns) -> pop_filter1_section*
$ws.Range("F10").Value = "pop_filter1_section2"
$ws.Range("G10").Value = "pop_filter1_section2_checkbox"
$ws.Range("H10").Value = "pop_filter1_section"

$ws.Range("F18").Value = "pop_filter1_section2"
$ws.Range("G18").Value = "pop_filter1_section2_checkbox"
$ws.Range("H18").Value = "pop_filter1_section"

# Update the view state: scroll position and selection
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F10").Select()
